# Auto-generated Excel COM-interop script applying scheduled-runner market/profit updates
# to the Behemoth_Profits workbook (per-sheet leve crafting tables).
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 8272.869000000001
$ws.Range("I86").Value = 8404.611000000001
$ws.Range("J86").Value = 7798.6
$ws.Range("K86").Value = 8404.611000000001
$ws.Range("L86").Value = 7798.6
$ws.Range("M86").Value = -7281.611000000001
$ws.Range("N86").Value = -10044.6
$ws.Range("H89").Value = 8272.869000000001
$ws.Range("I89").Value = 8404.611000000001
$ws.Range("J89").Value = 7798.6
$ws.Range("K89").Value = 42023.05500000001
$ws.Range("L89").Value = 38993
$ws.Range("M89").Value = -36407.05500000001
$ws.Range("N89").Value = -50225
$ws.Range("H99").Value = 1237
$ws.Range("I99").Value = 506.75
$ws.Range("J99").Value = 10000
$ws.Range("K99").Value = 1520.25
$ws.Range("L99").Value = 30000
$ws.Range("M99").Value = -22.25
$ws.Range("N99").Value = -32996
$ws.Range("H116").Value = 6992.533
$ws.Range("J116").Value = 7111
$ws.Range("L116").Value = 7111
$ws.Range("N116").Value = -13995
$ws.Range("H118").Value = 970.9
$ws.Range("I118").Value = 386
$ws.Range("J118").Value = 1555.8
$ws.Range("K118").Value = 1158
$ws.Range("L118").Value = 4667.4
$ws.Range("M118").Value = 499
$ws.Range("N118").Value = -7981.4
$ws.Range("H138").Value = 2618.4827
$ws.Range("I138").Value = 945.3
$ws.Range("J138").Value = 3499.1052
$ws.Range("K138").Value = 2835.9
$ws.Range("L138").Value = 10497.3156
$ws.Range("M138").Value = 2304.1
$ws.Range("N138").Value = -20777.3156

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2162.6428
$ws.Range("I45").Value = 2640.75
$ws.Range("K45").Value = 2640.75
$ws.Range("M45").Value = -2263.75
$ws.Range("H81").Value = 122000
$ws.Range("J81").Value = 122000
$ws.Range("L81").Value = 122000
$ws.Range("N81").Value = -123996
$ws.Range("H84").Value = 122000
$ws.Range("J84").Value = 122000
$ws.Range("L84").Value = 366000
$ws.Range("N84").Value = -375984
$ws.Range("H114").Value = 94132
$ws.Range("J114").Value = 94132
$ws.Range("L114").Value = 94132
$ws.Range("N114").Value = -102810
$ws.Range("H115").Value = 75209.5
$ws.Range("J115").Value = 75209.5
$ws.Range("L115").Value = 75209.5
$ws.Range("N115").Value = -78343.5
$ws.Range("H122").Value = 1285.25
$ws.Range("I122").Value = 1111.7142
$ws.Range("K122").Value = 3335.1426
$ws.Range("M122").Value = -885.1425999999997
$ws.Range("H124").Value = 37808.332
$ws.Range("J124").Value = 37808.332
$ws.Range("L124").Value = 37808.332
$ws.Range("N124").Value = -47628.332

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1834.2632
$ws.Range("I105").Value = 1060.8
$ws.Range("K105").Value = 1060.8
$ws.Range("M105").Value = 686.2
$ws.Range("H134").Value = 67937
$ws.Range("I134").Value = 2908.5
$ws.Range("K134").Value = 8725.5
$ws.Range("M134").Value = -6190.5

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2750
$ws.Range("I16").Value = 2500
$ws.Range("K16").Value = 2500
$ws.Range("M16").Value = -2213
$ws.Range("H44").Value = 30000
$ws.Range("I44").Value = 30000
$ws.Range("K44").Value = 30000
$ws.Range("M44").Value = -29558
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()
$ws.Range("H86").Value = 3000
$ws.Range("I86").Value = 3000
$ws.Range("J86").Value = 3000
$ws.Range("K86").Value = 3000
$ws.Range("L86").Value = 3000
$ws.Range("M86").Value = -1877
$ws.Range("N86").Value = -5246
$ws.Range("H89").Value = 3000
$ws.Range("I89").Value = 3000
$ws.Range("J89").Value = 3000
$ws.Range("K89").Value = 15000
$ws.Range("L89").Value = 15000
$ws.Range("M89").Value = -9384
$ws.Range("N89").Value = -26232
$ws.Range("H105").Value = 1841.5
$ws.Range("J105").Value = 1867.8572
$ws.Range("L105").Value = 1867.8572
$ws.Range("N105").Value = -5361.8572
$ws.Range("H113").Value = 2750
$ws.Range("I113").Value = 2500
$ws.Range("K113").Value = 2500
$ws.Range("M113").Value = -330
$ws.Range("H132").Value = 2950.8
$ws.Range("I132").Value = 2419.3684
$ws.Range("K132").Value = 7258.1052
$ws.Range("M132").Value = -4728.1052
$ws.Range("H134").Value = 315642.88
$ws.Range("I134").Value = 455693.88
$ws.Range("K134").Value = 1367081.64
$ws.Range("M134").Value = -1364546.64

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 528582
$ws.Range("I12").Value = 1679.1666
$ws.Range("K12").Value = 5037.4998
$ws.Range("M12").Value = -4864.4998
$ws.Range("H80").Value = 6000
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 6000
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 18000
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -19872
$ws.Range("H83").Value = 6000
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 6000
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 54000
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -63360
$ws.Range("H136").Value = 2283.1667
$ws.Range("I136").Value = 2283.1667
$ws.Range("K136").Value = 6849.500100000001
$ws.Range("M136").Value = -1749.500100000001
$ws.Range("H137").Value = 4719.467
$ws.Range("I137").Value = 6375.4287
$ws.Range("J137").Value = 3270.5
$ws.Range("K137").Value = 19126.2861
$ws.Range("L137").Value = 9811.5
$ws.Range("M137").Value = -14026.2861
$ws.Range("N137").Value = -20011.5
$ws.Range("H138").Value = 1539.5454
$ws.Range("I138").Value = 992.7778
$ws.Range("K138").Value = 2978.3334
$ws.Range("M138").Value = 2161.6666
$ws.Range("H139").Value = 99520.78
$ws.Range("I139").Value = 177137.4
$ws.Range("K139").Value = 531412.2
$ws.Range("M139").Value = -526272.2

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2867.9285
$ws.Range("I80").Value = 2906.75
$ws.Range("K80").Value = 2906.75
$ws.Range("M80").Value = -1908.75
$ws.Range("H83").Value = 2867.9285
$ws.Range("I83").Value = 2906.75
$ws.Range("K83").Value = 14533.75
$ws.Range("M83").Value = -9541.75
$ws.Range("H97").Value = 2389.8462
$ws.Range("I97").Value = 2172.5
$ws.Range("J97").Value = 4998
$ws.Range("K97").Value = 2172.5
$ws.Range("L97").Value = 4998
$ws.Range("M97").Value = -1676.5
$ws.Range("N97").Value = -5990
$ws.Range("H102").Value = 3985.8696
$ws.Range("I102").Value = 3061.9473
$ws.Range("J102").Value = 8374.5
$ws.Range("K102").Value = 3061.9473
$ws.Range("L102").Value = 8374.5
$ws.Range("M102").Value = -1439.9473
$ws.Range("N102").Value = -11618.5
$ws.Range("H132").Value = 45457876
$ws.Range("I132").Value = 62503532
$ws.Range("J132").Value = 2791.1667
$ws.Range("K132").Value = 187510596
$ws.Range("L132").Value = 8373.500100000001
$ws.Range("M132").Value = -187508066
$ws.Range("N132").Value = -13433.5001

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2684.0527
$ws.Range("I22").Value = 2730.5386
$ws.Range("K22").Value = 2730.5386
$ws.Range("M22").Value = -2435.5386
$ws.Range("H27").Value = 2684.0527
$ws.Range("I27").Value = 2730.5386
$ws.Range("K27").Value = 2730.5386
$ws.Range("M27").Value = -2623.5386
$ws.Range("H61").Value = 3400.2
$ws.Range("I61").Value = 3000.25
$ws.Range("K61").Value = 3000.25
$ws.Range("M61").Value = -2798.25
$ws.Range("H104").Value = 46730.145
$ws.Range("J104").Value = 46730.145
$ws.Range("L104").Value = 46730.145
$ws.Range("N104").Value = -53718.145
$ws.Range("H113").Value = 3400.2
$ws.Range("I113").Value = 3000.25
$ws.Range("K113").Value = 3000.25
$ws.Range("M113").Value = -830.25
$ws.Range("H122").Value = 4575.933
$ws.Range("I122").Value = 3895.3462
$ws.Range("J122").Value = 8999.75
$ws.Range("K122").Value = 11686.0386
$ws.Range("L122").Value = 26999.25
$ws.Range("M122").Value = -9236.0386
$ws.Range("N122").Value = -31899.25
$ws.Range("H136").Value = 96518.086
$ws.Range("I136").Value = 5321.2
$ws.Range("K136").Value = 15963.6
$ws.Range("M136").Value = -13413.6

